{"js": "// Apply the Mystic Spice Premium Chai Tea market-analysis report edits:\n// a handful of Chinese-text wording tweaks scattered through the body and\n// the benefits table. Each old string is unique in the document, so a\n// body-wide literal search + full-range replace is safe for every item.\n\nconst replacements = [\n  [\"\u4ecb\u7ecd\", \"\u7b80\u4ecb\"],\n  [\n    \"\u5b83\u8fd8\u627f\u8f7d\u7740\u6df1\u539a\u7684\u6587\u5316\u548c\u5386\u53f2\u610f\u4e49\uff0c\u5e38\u5e38\u4e0e\u597d\u5ba2\u3001\u53cb\u8c0a\u548c\u4f11\u95f2\u653e\u677e\u8054\u7cfb\u5728\u4e00\u8d77\u3002\",\n    \"\u5b83\u8fd8\u627f\u8f7d\u7740\u6df1\u539a\u7684\u6587\u5316\u548c\u5386\u53f2\u610f\u4e49\uff0c\u5e38\u4e0e\u597d\u5ba2\u3001\u53cb\u8c0a\u548c\u4f11\u95f2\u653e\u677e\u7d27\u5bc6\u76f8\u8fde\u3002\",\n  ],\n  [\"\u4ea7\u54c1\u63cf\u8ff0\", \"\u4ea7\u54c1\u8bf4\u660e\"],\n  [\n    \"\u5065\u5eb7\u589e\u5f3a\u6210\u5206\uff1a\u795e\u79d8\u9999\u6599\u8336\u4e2d\u7684\u6bcf\u79cd\u6210\u5206\u90fd\u7ecf\u8fc7\u7cbe\u5fc3\u6311\u9009\uff0c\u4ee5\u4fbf\u53d1\u6325\u5176\u5929\u7136\u7684\u5065\u5eb7\u76ca\u5904\u3002\",\n    \"\u6210\u5206\u66f4\u52a0\u5065\u5eb7\uff1a\u795e\u79d8\u9999\u6599\u5976\u8336\u81fb\u9009\u81ea\u7136\u539f\u6599\uff0c\u6709\u5229\u4e8e\u5065\u5eb7\u3002\",\n  ],\n  [\n    \"\u6d53\u90c1\u7684\u9999\u6c14\u548c\u98ce\u5473\uff1a\u6211\u4eec\u7684\u8336\u53e3\u611f\u6e29\u548c\uff0c\u9999\u4e2d\u5e26\u8fa3\uff0c\u5473\u9053\u6d53\u90c1\uff0c\u4ee4\u4eba\u632f\u594b\uff0c\u662f\u5f00\u542f\u4e00\u5929\u6216\u665a\u4e0a\u653e\u677e\u7684\u5b8c\u7f8e\u996e\u6599\u3002\",\n    \"\u9999\u6c14\u6d53\u90c1\u3001\u53e3\u5473\u9187\u539a\uff1a\u6211\u4eec\u7684\u5976\u8336\u6c14\u5473\u6e29\u8f9b\u3001\u53e3\u611f\u9187\u539a\uff0c\u63d0\u795e\u9192\u8111\uff0c\u662f\u5f00\u542f\u7f8e\u597d\u4e00\u5929\u6216\u665a\u4e0a\u653e\u677e\u8eab\u5fc3\u7684\u5b8c\u7f8e\u996e\u54c1\u3002\",\n  ],\n  [\n    \"\u591a\u79cd\u51b2\u6ce1\u9009\u9879\uff1a\u65e0\u8bba\u559c\u6b22\u70ed\u6c14\u817e\u817e\u7684\u8336\u3001\u6e05\u723d\u7684\u51b0\u8336\u8fd8\u662f\u5976\u6cb9\u62ff\u94c1\uff0c\u6211\u4eec\u7684\u6df7\u5408\u54c1\u90fd\u53ef\u7075\u6d3b\u642d\u914d\uff0c\u53ef\u4ee5\u6ee1\u8db3\u4efb\u4f55\u559c\u597d\u3002\",\n    \"\u591a\u5143\u5316\u7684\u70f9\u5236\u9009\u9879\uff1a\u65e0\u8bba\u4f60\u662f\u559c\u6b22\u6e29\u70ed\u7684\u5976\u8336\uff0c\u8fd8\u662f\u4ee4\u4eba\u8033\u76ee\u4e00\u65b0\u7684\u51b0\u8336\uff0c\u6216\u8005\u662f\u5976\u6cb9\u62ff\u94c1\uff0c\u8fd9\u6b3e\u4ea7\u54c1\u53ef\u4ee5\u6ee1\u8db3\u4efb\u4f55\u504f\u597d\u3002\",\n  ],\n  [\n    \"\u53ef\u6301\u7eed\u91c7\u8d2d\uff1a\u6211\u4eec\u81f4\u529b\u4e8e\u53ef\u6301\u7eed\u53d1\u5c55\uff0c\u4ece\u4ece\u4e8b\u6709\u673a\u519c\u4e1a\u7684\u5c0f\u578b\u519c\u573a\u91c7\u8d2d\u539f\u6599\uff0c\u4e0d\u4ec5\u53ef\u786e\u4fdd\u6700\u4f73\u54c1\u8d28\uff0c\u800c\u4e14\u53ef\u9020\u798f\u5168\u4eba\u7c7b\u3002\",\n    \"\u539f\u6599\u53ef\u6301\u7eed\uff1a\u6211\u4eec\u6ce8\u91cd\u53ef\u6301\u7eed\u6027\uff0c\u4ece\u5c0f\u578b\u519c\u573a\u91c7\u8d2d\u539f\u6599\uff0c\u575a\u6301\u6709\u673a\u519c\u4e1a\uff0c\u4e0d\u4ec5\u80fd\u591f\u786e\u4fdd\u6781\u4f73\u54c1\u8d28\uff0c\u800c\u4e14\u53ef\u4ee5\u786e\u4fdd\u5bf9\u6211\u4eec\u7684\u661f\u7403\u6709\u76ca\u3002\",\n  ],\n  [\n    \"\u4f18\u96c5\u7684\u5305\u88c5\uff1a\u795e\u79d8\u9999\u6599\u8336\u7684\u5305\u88c5\u8bbe\u8ba1\u7cbe\u7f8e\u3001\u73af\u4fdd\uff0c\u662f\u8d60\u9001\u7231\u8336\u4eba\u58eb\u7684\u7406\u60f3\u793c\u7269\u4e4b\u9009\uff0c\u4e5f\u662f\u81ea\u884c\u4eab\u7528\u7684\u5962\u4f88\u9009\u62e9\u3002\",\n    \"\u5305\u88c5\u7cbe\u81f4\uff1a\u795e\u79d8\u9999\u6599\u5370\u5ea6\u5976\u8336\u8bbe\u8ba1\u7cbe\u7f8e\uff0c\u91c7\u7528\u751f\u6001\u53cb\u597d\u7684\u5305\u88c5\u65b9\u5f0f\uff0c\u56e0\u6b64\u662f\u9001\u7ed9\u8336\u53f6\u7231\u597d\u8005\u7684\u7406\u60f3\u793c\u7269\uff0c\u4e5f\u662f\u9001\u7ed9\u81ea\u5df1\u7684\u5962\u534e\u4e4b\u9009\u3002\",\n  ],\n  [\n    \"\u5ba2\u6237\u6ee1\u610f\u5ea6\u4fdd\u8bc1\uff1a\u6211\u4eec\u652f\u6301\u8fd9\u6b3e\u4ea7\u54c1\uff0c\u5e76\u63d0\u4f9b\u6ee1\u610f\u5ea6\u4fdd\u8bc1\u3002\",\n    \"\u5ba2\u6237\u6ee1\u610f\u5ea6\u4fdd\u8bc1\uff1a\u6211\u4eec\u4e3a\u4ea7\u54c1\u80cc\u4e66\uff0c\u63d0\u4f9b\u4ee4\u4eba\u6ee1\u610f\u7684\u4fdd\u8bc1\u3002\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the Mystic Spice Premium Chai Tea market-analysis report edits:\n# a handful of Chinese-text wording tweaks scattered through the body and\n# the benefits table. Each old string is unique in the document, so a\n# document-wide Find/Replace (wdReplaceAll) is safe for every item.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"\u4ecb\u7ecd\", \"\u7b80\u4ecb\"),\n    @(\"\u5b83\u8fd8\u627f\u8f7d\u7740\u6df1\u539a\u7684\u6587\u5316\u548c\u5386\u53f2\u610f\u4e49\uff0c\u5e38\u5e38\u4e0e\u597d\u5ba2\u3001\u53cb\u8c0a\u548c\u4f11\u95f2\u653e\u677e\u8054\u7cfb\u5728\u4e00\u8d77\u3002\", \"\u5b83\u8fd8\u627f\u8f7d\u7740\u6df1\u539a\u7684\u6587\u5316\u548c\u5386\u53f2\u610f\u4e49\uff0c\u5e38\u4e0e\u597d\u5ba2\u3001\u53cb\u8c0a\u548c\u4f11\u95f2\u653e\u677e\u7d27\u5bc6\u76f8\u8fde\u3002\"),\n    @(\"\u4ea7\u54c1\u63cf\u8ff0\", \"\u4ea7\u54c1\u8bf4\u660e\"),\n    @(\"\u5065\u5eb7\u589e\u5f3a\u6210\u5206\uff1a\u795e\u79d8\u9999\u6599\u8336\u4e2d\u7684\u6bcf\u79cd\u6210\u5206\u90fd\u7ecf\u8fc7\u7cbe\u5fc3\u6311\u9009\uff0c\u4ee5\u4fbf\u53d1\u6325\u5176\u5929\u7136\u7684\u5065\u5eb7\u76ca\u5904\u3002\", \"\u6210\u5206\u66f4\u52a0\u5065\u5eb7\uff1a\u795e\u79d8\u9999\u6599\u5976\u8336\u81fb\u9009\u81ea\u7136\u539f\u6599\uff0c\u6709\u5229\u4e8e\u5065\u5eb7\u3002\"),\n    @(\"\u6d53\u90c1\u7684\u9999\u6c14\u548c\u98ce\u5473\uff1a\u6211\u4eec\u7684\u8336\u53e3\u611f\u6e29\u548c\uff0c\u9999\u4e2d\u5e26\u8fa3\uff0c\u5473\u9053\u6d53\u90c1\uff0c\u4ee4\u4eba\u632f\u594b\uff0c\u662f\u5f00\u542f\u4e00\u5929\u6216\u665a\u4e0a\u653e\u677e\u7684\u5b8c\u7f8e\u996e\u6599\u3002\", \"\u9999\u6c14\u6d53\u90c1\u3001\u53e3\u5473\u9187\u539a\uff1a\u6211\u4eec\u7684\u5976\u8336\u6c14\u5473\u6e29\u8f9b\u3001\u53e3\u611f\u9187\u539a\uff0c\u63d0\u795e\u9192\u8111\uff0c\u662f\u5f00\u542f\u7f8e\u597d\u4e00\u5929\u6216\u665a\u4e0a\u653e\u677e\u8eab\u5fc3\u7684\u5b8c\u7f8e\u996e\u54c1\u3002\"),\n    @(\"\u591a\u79cd\u51b2\u6ce1\u9009\u9879\uff1a\u65e0\u8bba\u559c\u6b22\u70ed\u6c14\u817e\u817e\u7684\u8336\u3001\u6e05\u723d\u7684\u51b0\u8336\u8fd8\u662f\u5976\u6cb9\u62ff\u94c1\uff0c\u6211\u4eec\u7684\u6df7\u5408\u54c1\u90fd\u53ef\u7075\u6d3b\u642d\u914d\uff0c\u53ef\u4ee5\u6ee1\u8db3\u4efb\u4f55\u559c\u597d\u3002\", \"\u591a\u5143\u5316\u7684\u70f9\u5236\u9009\u9879\uff1a\u65e0\u8bba\u4f60\u662f\u559c\u6b22\u6e29\u70ed\u7684\u5976\u8336\uff0c\u8fd8\u662f\u4ee4\u4eba\u8033\u76ee\u4e00\u65b0\u7684\u51b0\u8336\uff0c\u6216\u8005\u662f\u5976\u6cb9\u62ff\u94c1\uff0c\u8fd9\u6b3e\u4ea7\u54c1\u53ef\u4ee5\u6ee1\u8db3\u4efb\u4f55\u504f\u597d\u3002\"),\n    @(\"\u53ef\u6301\u7eed\u91c7\u8d2d\uff1a\u6211\u4eec\u81f4\u529b\u4e8e\u53ef\u6301\u7eed\u53d1\u5c55\uff0c\u4ece\u4ece\u4e8b\u6709\u673a\u519c\u4e1a\u7684\u5c0f\u578b\u519c\u573a\u91c7\u8d2d\u539f\u6599\uff0c\u4e0d\u4ec5\u53ef\u786e\u4fdd\u6700\u4f73\u54c1\u8d28\uff0c\u800c\u4e14\u53ef\u9020\u798f\u5168\u4eba\u7c7b\u3002\", \"\u539f\u6599\u53ef\u6301\u7eed\uff1a\u6211\u4eec\u6ce8\u91cd\u53ef\u6301\u7eed\u6027\uff0c\u4ece\u5c0f\u578b\u519c\u573a\u91c7\u8d2d\u539f\u6599\uff0c\u575a\u6301\u6709\u673a\u519c\u4e1a\uff0c\u4e0d\u4ec5\u80fd\u591f\u786e\u4fdd\u6781\u4f73\u54c1\u8d28\uff0c\u800c\u4e14\u53ef\u4ee5\u786e\u4fdd\u5bf9\u6211\u4eec\u7684\u661f\u7403\u6709\u76ca\u3002\"),\n    @(\"\u4f18\u96c5\u7684\u5305\u88c5\uff1a\u795e\u79d8\u9999\u6599\u8336\u7684\u5305\u88c5\u8bbe\u8ba1\u7cbe\u7f8e\u3001\u73af\u4fdd\uff0c\u662f\u8d60\u9001\u7231\u8336\u4eba\u58eb\u7684\u7406\u60f3\u793c\u7269\u4e4b\u9009\uff0c\u4e5f\u662f\u81ea\u884c\u4eab\u7528\u7684\u5962\u4f88\u9009\u62e9\u3002\", \"\u5305\u88c5\u7cbe\u81f4\uff1a\u795e\u79d8\u9999\u6599\u5370\u5ea6\u5976\u8336\u8bbe\u8ba1\u7cbe\u7f8e\uff0c\u91c7\u7528\u751f\u6001\u53cb\u597d\u7684\u5305\u88c5\u65b9\u5f0f\uff0c\u56e0\u6b64\u662f\u9001\u7ed9\u8336\u53f6\u7231\u597d\u8005\u7684\u7406\u60f3\u793c\u7269\uff0c\u4e5f\u662f\u9001\u7ed9\u81ea\u5df1\u7684\u5962\u534e\u4e4b\u9009\u3002\"),\n    @(\"\u5ba2\u6237\u6ee1\u610f\u5ea6\u4fdd\u8bc1\uff1a\u6211\u4eec\u652f\u6301\u8fd9\u6b3e\u4ea7\u54c1\uff0c\u5e76\u63d0\u4f9b\u6ee1\u610f\u5ea6\u4fdd\u8bc1\u3002\", \"\u5ba2\u6237\u6ee1\u610f\u5ea6\u4fdd\u8bc1\uff1a\u6211\u4eec\u4e3a\u4ea7\u54c1\u80cc\u4e66\uff0c\u63d0\u4f9b\u4ee4\u4eba\u6ee1\u610f\u7684\u4fdd\u8bc1\u3002\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
